# "Examples from 11 to 19 added"
# Adds three new columns (VQ_best, comp_offs_y, comp_offs_x) with data for
# rows 2-4, widens the used range accordingly, and swaps the values that used
# to live in rows 3 and 4 (columns M:R) so that what was row 4 becomes row 3
# and vice versa.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column headers (T1:V1) ---------------------------------------
$ws.Range("T1").Value = "VQ_best"
$ws.Range("U1").Value = "comp_offs_y"
$ws.Range("V1").Value = "comp_offs_x"

# --- New column data -----------------------------------------------------
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 0

$ws.Range("U3").Value = 125
$ws.Range("V3").Value = 167

$ws.Range("U4").Value = 120
$ws.Range("V4").Value = 167

# --- Swap the candidate-motion-vector rows (M3:R3 <-> M4:R4) -------------
# Row 3 used to hold what is now row 4's data, and vice versa.
$cols = @("M", "N", "O", "P", "Q", "R")
$newRow3 = @(-48, 112, 0, -96, -144, 128)
$newRow4 = @(-64, 128, -128, 64, -128, 144)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $newRow3[$i]
    $ws.Range($cols[$i] + "4").Value = $newRow4[$i]
}

# --- Column widths for the two newly introduced columns ------------------
$ws.Range("U1").ColumnWidth = 14
$ws.Range("V1").ColumnWidth = 12.333333333333332

# --- Selection / scroll position, matching the saved view ----------------
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("V6").Select()
